$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '76.448.04'
$ws.Range("E2").Value = '  +0.20%  '

# Row 3
$ws.Range("D3").Value = '2.919.29'
$ws.Range("E3").Value = '  +0.64%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '197.71'
$ws.Range("E5").Value = '  +0.49%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '593.31'
$ws.Range("E6").Value = '  -1.33%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("E8").Value = '  -1.49%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.198'
$ws.Range("E9").Value = '  +2.47%  '

# Row 10
$ws.Range("D10").Value = '2.918.62'
$ws.Range("E10").Value = '  +0.65%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("E11").Value = '  +10.00%  '

# Row 12
$ws.Range("E12").Value = '  +0.24%  '

# Row 13
$ws.Range("D13").Value = '3.459.94'
$ws.Range("E13").Value = '  +1.62%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.85'
$ws.Range("E14").Value = '  -2.22%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '76.481.12'
$ws.Range("E15").Value = '  +0.37%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.12'
$ws.Range("E16").Value = '  +1.75%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000188'
$ws.Range("E17").Value = '  -1.70%  '

# Row 18
$ws.Range("D18").Value = '2.924.15'
$ws.Range("E18").Value = '  +0.72%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +5.73%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.70'
$ws.Range("E20").Value = '  -3.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.36'
$ws.Range("E21").Value = '  -3.33%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  +2.55%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.24'
$ws.Range("E23").Value = '  -3.00%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.85'
$ws.Range("E24").Value = '  -0.20%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.09%  '

# Row 26
$ws.Range("D26").Value = '3.080.89'
$ws.Range("E26").Value = '  +1.59%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.23'
$ws.Range("E27").Value = '  -1.37%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.56'
$ws.Range("E28").Value = '  -2.58%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000106'
$ws.Range("E29").Value = '  -2.38%  '

# Row 30
$ws.Range("E30").Value = '  +0.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.26'
$ws.Range("E31").Value = '  +5.24%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.36'
$ws.Range("E32").Value = '  -4.21%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '500.22'
$ws.Range("E33").Value = '  -2.94%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -0.60%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.01%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.43'
$ws.Range("E36").Value = '  +0.03%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.04'
$ws.Range("E37").Value = '  -1.12%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.110'
$ws.Range("E38").Value = '  +17.99%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.390'
$ws.Range("E39").Value = '  +11.43%  '

# Row 40
$ws.Range("E40").Value = '  +1.29%  '

# Row 41
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.110'
$ws.Range("E42").Value = '  -6.14%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '178.34'
$ws.Range("E43").Value = '  -3.22%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.89'
$ws.Range("E44").Value = '  -4.29%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.64'
$ws.Range("E45").Value = '  -3.16%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.98'
$ws.Range("E46").Value = '  -0.50%  '

# Row 47
$ws.Range("E47").Value = '  -5.14%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.587'
$ws.Range("E48").Value = '  +0.21%  '

# Row 49
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.30'
$ws.Range("E49").Value = '  -4.03%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.84'
$ws.Range("E50").Value = '  +1.49%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.18'
$ws.Range("E51").Value = '  +2.49%  '
